$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns I (9) through N (14) on row 1.
# I1/J1 are brand-new "rate" columns, K1 keeps the old CHEQUE.NUMBER
# header (shifted right from I1), and L1:N1 are three more new columns.
$ws.Range("I1").Value = "TREASURY.RATE"
$ws.Range("J1").Value = "CUST.RATE"
$ws.Range("K1").Value = "CHEQUE.NUMBER"
$ws.Range("L1").Value = "DD.ADDRESS:1"
$ws.Range("M1").Value = "PURP.REMITT:1"
$ws.Range("N1").Value = "REL.BENEFICIARY"

# The cheque number data value moves from I2 to K2 to line up with the
# relocated CHEQUE.NUMBER header; I2 becomes blank.
$ws.Range("K2").Value = $ws.Range("I2").Value2
$ws.Range("I2").ClearContents()

# Give the two brand new columns (I:J) a manual width matching column H's
# width, mirroring the customWidth (non-bestFit) columns added alongside
# the new headers.
$ws.Range("I1:J1").ColumnWidth = 15.83

# Reflect the author's last UI selection when the sheet was saved.
$ws.Range("F18").Select() | Out-Null
